$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.26419815983392
$ws.Range("D2").Value = 0.531475020769009
$ws.Range("B3").Value = 24253.1967591745
$ws.Range("B4").Value = 545.931539920301
$ws.Range("D4").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000106375289940729
